$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing data rows 2-52 down to 3-53)
$ws.Rows("2:2").Insert()

# Reset formatting on the freshly inserted row, then copy the date-column format
# from row 3 (column A) onto the new A2 cell so it matches the rest of column A.
# E2 has no data in the new row, so fully clear it (removes the cell, matching the
# sparse layout used by every other row that lacks a y_1_forecast value).
$ws.Range("A2:D2").ClearFormats()
$ws.Range("E2").Clear()
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 2 (date serial 2007-06-24 / y_0=2007 / y_0_forecast / y_1=2008)
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 2.070003986395053
$ws.Range("D2").Value = 2008

# Refresh the recomputed y_0_forecast (C) / y_1_forecast (E) values for every shifted row
$ws.Range("C4").Value = 0.517569958955022
$ws.Range("C6").Value = -3.956152295564896
$ws.Range("C7").Value = 0.7825601129312298
$ws.Range("E7").Value = -0.2328395085068102
$ws.Range("C8").Value = 1.234995474941392
$ws.Range("E8").Value = 0.7327527981400461
$ws.Range("C9").Value = 1.171373351779592
$ws.Range("E9").Value = 0.6458049593451864
$ws.Range("C10").Value = 0.899360810820804
$ws.Range("E10").Value = 0.4113848771853501
$ws.Range("C11").Value = 0.4797371259343874
$ws.Range("E11").Value = 0.5448268972618964
$ws.Range("C12").Value = 0.9010266119894084
$ws.Range("E12").Value = 1.203907967581297
$ws.Range("C13").Value = 0.1088602047940146
$ws.Range("E13").Value = 0.1671551101610103
$ws.Range("C14").Value = 0.02019328874804938
$ws.Range("E14").Value = -0.03860754389363175
$ws.Range("C15").Value = -0.075394216261504
$ws.Range("E15").Value = -0.06463514052835739
$ws.Range("C16").Value = 0.1729981757035093
$ws.Range("E16").Value = 0.2629870913912535
$ws.Range("C17").Value = -0.07548837955325682
$ws.Range("E17").Value = 0.05928147027902675
$ws.Range("C18").Value = 0.09752710595589686
$ws.Range("E18").Value = 0.1987429576382649
$ws.Range("C19").Value = -0.5735475396625112
$ws.Range("E19").Value = 0.02406984837131088
$ws.Range("C20").Value = -0.5280591151586633
$ws.Range("E20").Value = 0.1903092973221776
$ws.Range("C21").Value = 0.02883110668334687
$ws.Range("E21").Value = 0.241498802789164
$ws.Range("C22").Value = 0.07201851318385799
$ws.Range("E22").Value = 0.2843016498274009
$ws.Range("C23").Value = 0.8258453722611359
$ws.Range("E23").Value = 0.288064297781454
$ws.Range("C24").Value = 0.7252300059688022
$ws.Range("E24").Value = 0.2052430644269299
$ws.Range("C25").Value = 0.4640111827386662
$ws.Range("E25").Value = -0.1620498231152179
$ws.Range("C26").Value = 0.3727661260635617
$ws.Range("E26").Value = -0.8612142616933327
$ws.Range("C27").Value = -0.9065026814729205
$ws.Range("E27").Value = -0.1365403697986656
$ws.Range("C28").Value = -0.6243248145489155
$ws.Range("E28").Value = 0.07475225043114264
$ws.Range("C29").Value = -0.7901161779547028
$ws.Range("E29").Value = -0.1949185821441768
$ws.Range("C30").Value = -0.801759526476209
$ws.Range("E30").Value = 0.06491682578968483
$ws.Range("C31").Value = -0.2930109800340586
$ws.Range("E31").Value = -0.1644047560850792
$ws.Range("C32").Value = -0.6491730431770759
$ws.Range("E32").Value = -0.4370777949570193
$ws.Range("C33").Value = -1.503583188367719
$ws.Range("E33").Value = -0.5877299932847579
$ws.Range("C34").Value = -1.103489789942047
$ws.Range("E34").Value = 1.323658311025055
$ws.Range("C35").Value = 2.632698787096288
$ws.Range("E35").Value = 0.4893323826990148
$ws.Range("C36").Value = 1.311904119834839
$ws.Range("E36").Value = -0.1159018519404809
$ws.Range("C37").Value = 1.067534122491809
$ws.Range("E37").Value = -0.112565850764601
$ws.Range("C38").Value = 0.9704846793491928
$ws.Range("E38").Value = -0.6989646400249128
$ws.Range("C39").Value = -0.3568974718008655
$ws.Range("E39").Value = 0.04689880979749095
$ws.Range("C40").Value = -1.338216592160768
$ws.Range("E40").Value = -0.2487719682984557
$ws.Range("C41").Value = -0.9795431199870586
$ws.Range("E41").Value = -0.3324688493351879
$ws.Range("C42").Value = -0.7009264669202708
$ws.Range("E42").Value = -0.0234350458557242
$ws.Range("C43").Value = -0.05370673382950608
$ws.Range("E43").Value = -0.03417477517112522
$ws.Range("C44").Value = 0.2512652100014945
$ws.Range("E44").Value = 0.041441321352087
$ws.Range("C45").Value = 0.1829021030556488
$ws.Range("E45").Value = 0.01245506629512505
$ws.Range("C46").Value = 0.3928252664241905
$ws.Range("E46").Value = 0.196134499498668
$ws.Range("C47").Value = 0.2464401331885524
$ws.Range("E47").Value = -0.02784774425726999
$ws.Range("C48").Value = 0.9693451788297391
$ws.Range("E48").Value = 0.08029846083614789
$ws.Range("C49").Value = 0.6979546684258597
$ws.Range("E49").Value = -0.04101132194430646
$ws.Range("C50").Value = 0.3224026462283813
$ws.Range("E50").Value = -0.6671574593505647
$ws.Range("C51").Value = -2.205730080079726
$ws.Range("E51").Value = -0.2015885781823656
$ws.Range("C52").Value = -1.551451534890558
$ws.Range("E52").Value = -0.2835476113072333
$ws.Range("C53").Value = -2.11737366557071
$ws.Range("E53").Value = -0.2835476113072333
